# Daily attendance processing - 2025-10-07 23:39:47
# Applies refreshed attendance-sync values to the session analysis sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / numeric cell updates -------------------------------------
$ws.Range("G2").Value = "system, backup@backdoor.com, System"
$ws.Range("G3").Value = "dnasr281@gmail.com, System"
$ws.Range("H3").Value = "53/53"
$ws.Range("G4").Value = "backup@backdoor.com, System"
$ws.Range("H4").Value = "40/53"
$ws.Range("G5").Value = "backup@backdoor.com, System"
$ws.Range("H6").Value = "45/53"
$ws.Range("L6").Value = 83
$ws.Range("L7").Value = 1
$ws.Range("G11").Value = "dnasr281@gmail.com, System"
$ws.Range("G14").Value = "dnasr281@gmail.com, System"
$ws.Range("G15").Value = "dnasr281@gmail.com, System"
$ws.Range("O18").Value = 13
$ws.Range("P18").Value = 0
$ws.Range("O19").Value = 13
$ws.Range("P19").Value = 0
$ws.Range("O20").Value = 13
$ws.Range("P20").Value = 0
$ws.Range("G29").Value = "system, backup@backdoor.com, System"
$ws.Range("H29").Value = "35/56"
$ws.Range("G30").Value = "dnasr281@gmail.com, System"
$ws.Range("H30").Value = "42/56"
$ws.Range("G32").Value = "backup@backdoor.com, System"
$ws.Range("H32").Value = "38/56"
$ws.Range("H33").Value = "43/56"
$ws.Range("G38").Value = "dnasr281@gmail.com, System"
$ws.Range("H38").Value = "32/56"
$ws.Range("H39").Value = "36/56"
$ws.Range("G41").Value = "dnasr281@gmail.com, System"
$ws.Range("H41").Value = "44/56"
$ws.Range("G42").Value = "dnasr281@gmail.com, System"
$ws.Range("H42").Value = "45/56"
$ws.Range("G56").Value = "system, backup@backdoor.com, System"
$ws.Range("H56").Value = "55/55"
$ws.Range("G57").Value = "dnasr281@gmail.com, System"
$ws.Range("G58").Value = "backup@backdoor.com, System"
$ws.Range("G59").Value = "backup@backdoor.com, System"
$ws.Range("H60").Value = "36/55"
$ws.Range("G65").Value = "dnasr281@gmail.com, System"
$ws.Range("G68").Value = "dnasr281@gmail.com, System"
$ws.Range("G69").Value = "dnasr281@gmail.com, System"
$ws.Range("G84").Value = "backup@backdoor.com, System"
$ws.Range("H84").Value = "56/56"
$ws.Range("G85").Value = "backup@backdoor.com, System"
$ws.Range("G86").Value = "dnasr281@gmail.com, System"
$ws.Range("H86").Value = "37/56"
$ws.Range("G89").Value = "dnasr281@gmail.com, System"
$ws.Range("G90").Value = "dnasr281@gmail.com, admin@admin.com"
$ws.Range("G93").Value = "dnasr281@gmail.com, System"
$ws.Range("G110").Value = "backup@backdoor.com, System"
$ws.Range("G111").Value = "backup@backdoor.com, System"
$ws.Range("H111").Value = "55/55"
$ws.Range("G112").Value = "dnasr281@gmail.com, System"
$ws.Range("G115").Value = "dnasr281@gmail.com, System"
$ws.Range("G116").Value = "dnasr281@gmail.com, admin@admin.com"
$ws.Range("G119").Value = "dnasr281@gmail.com, System"
$ws.Range("H119").Value = "41/55"
$ws.Range("G136").Value = "backup@backdoor.com, System"
$ws.Range("G137").Value = "backup@backdoor.com, System"
$ws.Range("G138").Value = "dnasr281@gmail.com, System"
$ws.Range("G141").Value = "dnasr281@gmail.com, System"
$ws.Range("G142").Value = "dnasr281@gmail.com, admin@admin.com"
$ws.Range("G145").Value = "dnasr281@gmail.com, System"

# --- Percentage values stored as literal text (match source formatting) ---
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "52.2%"
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "67.1%"
$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = "69.7%"
$ws.Range("S16").NumberFormat = "@"
$ws.Range("S16").Value = "65.1%"
$ws.Range("S17").NumberFormat = "@"
$ws.Range("S17").Value = "61.5%"
$ws.Range("R18").NumberFormat = "@"
$ws.Range("R18").Value = "50.0%"
$ws.Range("S18").NumberFormat = "@"
$ws.Range("S18").Value = "69.1%"
$ws.Range("R19").NumberFormat = "@"
$ws.Range("R19").Value = "50.0%"
$ws.Range("S19").NumberFormat = "@"
$ws.Range("S19").Value = "70.3%"
$ws.Range("R20").NumberFormat = "@"
$ws.Range("R20").Value = "50.0%"
$ws.Range("S20").NumberFormat = "@"
$ws.Range("S20").Value = "67.5%"

# --- Rows newly recorded by the System sync (A:I formatting + values) ------
$ws.Range("A2:I2").Copy()
$ws.Range("A95:I95").PasteSpecial(-4122)
$ws.Range("G95").Value = "System"
$ws.Range("H95").Value = "1/56"
$ws.Range("I95").Value = "Recorded"

$ws.Range("A2:I2").Copy()
$ws.Range("A121:I121").PasteSpecial(-4122)
$ws.Range("G121").Value = "System"
$ws.Range("H121").Value = "2/55"
$ws.Range("I121").Value = "Recorded"

$ws.Range("A2:I2").Copy()
$ws.Range("A147:I147").PasteSpecial(-4122)
$ws.Range("G147").Value = "System"
$ws.Range("H147").Value = "2/57"
$ws.Range("I147").Value = "Recorded"

